$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Bdnf", "Ngfr", "FAPs", 1, 0.3333333333333333, 0.03885866666666667, 0.116576, 0.09340097618505853, 0.09340097618505853, 3, 1, 3.362744666666666, 10.088234, 0.7488888671136141, 0.748888867113614, 0.1306717740871111, 1.176045966784, 0.06994695124253414, 0.06994695124253412),
    @("ECs", "Bdnf", "Ngfr", "MuSCs", 1, 0.3333333333333333, 0.03885866666666667, 0.116576, 0.09340097618505853, 0.09340097618505853, 3, 1, 1.127567333333333, 3.382702, 0.251111132886386, 0.2511111328863859, 0.04381576315022222, 0.394341868352, 0.02345402494252441, 0.0234540249425244),
    @("MuSCs", "Bdnf", "Ngfr", "FAPs", 2, 0.6666666666666666, 0.3771826666666667, 1.131548, 0.9065990238149415, 0.9065990238149415, 3, 1, 3.362744666666666, 10.088234, 0.7488888671136141, 0.748888867113614, 1.268369000692444, 11.415321006232, 0.67894191587108, 0.6789419158710799),
    @("MuSCs", "Bdnf", "Ngfr", "MuSCs", 2, 0.6666666666666666, 0.3771826666666667, 1.131548, 0.9065990238149415, 0.9065990238149415, 3, 1, 1.127567333333333, 3.382702, 0.251111132886386, 0.2511111328863859, 0.4252988536328889, 3.827689682696, 0.2276571079438616, 0.2276571079438615)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
